$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be stored as text so values like "1.008" are not
# auto-converted to numbers by Excel (they must remain literal strings, matching
# the original inline-string cell type).
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

# Apply the updated cell values (Coin / Link / Price / Volume(1h)).
$ws.Range('D2').Value = '21.848.62'
$ws.Range('E2').Value = '  -1.17%  '
$ws.Range('D3').Value = '1.551.61'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +0.79%  '
$ws.Range('D5').Value = '1.004'
$ws.Range('E5').Value = '  +0.42%  '
$ws.Range('D6').Value = '288.35'
$ws.Range('E6').Value = '  +0.23%  '
$ws.Range('D7').Value = '0.3887'
$ws.Range('E7').Value = '  +1.98%  '
$ws.Range('D8').Value = '0.3184'
$ws.Range('E8').Value = '  -2.98%  '
$ws.Range('D9').Value = '43.64'
$ws.Range('E9').Value = '  +1.05%  '
$ws.Range('D10').Value = '0.07122'
$ws.Range('E10').Value = '  -2.90%  '
$ws.Range('D11').Value = '1.054'
$ws.Range('E11').Value = '  -7.04%  '
$ws.Range('D12').Value = '1.007'
$ws.Range('E12').Value = '  +0.78%  '
$ws.Range('D13').Value = '5.585'
$ws.Range('E13').Value = '  -3.67%  '
$ws.Range('D14').Value = '18.43'
$ws.Range('E14').Value = '  -8.18%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '6.591'
$ws.Range('E15').Value = '  -2.88%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.543.31'
$ws.Range('E16').Value = '  +0.60%  '
$ws.Range('D17').Value = '0.00001094'
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').Value = '0.06548'
$ws.Range('E18').Value = '  -0.80%  '
$ws.Range('D19').Value = '82.67'
$ws.Range('E19').Value = '  -3.50%  '
$ws.Range('E20').Value = '  +0.35%  '
$ws.Range('D21').Value = '6.142'
$ws.Range('E21').Value = '  -3.89%  '
$ws.Range('D22').Value = '15.19'
$ws.Range('E22').Value = '  -5.29%  '
$ws.Range('D23').Value = '10.97'
$ws.Range('E23').Value = '  -5.91%  '
$ws.Range('D24').Value = '2.373'
$ws.Range('E24').Value = '  +2.44%  '
$ws.Range('D25').Value = '21.881.32'
$ws.Range('E25').Value = '  -0.98%  '
$ws.Range('D26').Value = '2.353'
$ws.Range('E26').Value = '  -6.37%  '
$ws.Range('D27').Value = '145.29'
$ws.Range('E27').Value = '  -3.25%  '
$ws.Range('D28').Value = '18.34'
$ws.Range('E28').Value = '  -3.78%  '
$ws.Range('D29').Value = '4.854'
$ws.Range('E29').Value = '  -1.10%  '
$ws.Range('D30').Value = '1.717.78'
$ws.Range('E30').Value = '  +0.35%  '
$ws.Range('D31').Value = '116.89'
$ws.Range('E31').Value = '  -3.61%  '
$ws.Range('D32').Value = '0.9672'
$ws.Range('E32').Value = '  -9.11%  '
$ws.Range('D33').Value = '5.802'
$ws.Range('E33').Value = '  -1.62%  '
$ws.Range('D34').Value = '0.08173'
$ws.Range('E34').Value = '  -0.49%  '
$ws.Range('B35').Value = 'FraxShare'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D35').Value = '8.966'
$ws.Range('E35').Value = '  -3.25%  '
$ws.Range('B36').Value = 'WEMIXTOKEN'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = '1.595'
$ws.Range('E36').Value = '  -14.53%  '
$ws.Range('D37').Value = '0.02216'
$ws.Range('E37').Value = '  -4.50%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = '5.048'
$ws.Range('E38').Value = '  -4.24%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.05961'
$ws.Range('E39').Value = '  -4.15%  '
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').Value = '0.2015'
$ws.Range('E40').Value = '  -6.59%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '1.178'
$ws.Range('E41').Value = '  -5.69%  '
$ws.Range('D42').Value = '1.003'
$ws.Range('E42').Value = '  +0.32%  '
$ws.Range('D43').Value = '10.44'
$ws.Range('E43').Value = '  -5.02%  '
$ws.Range('D44').Value = '0.5706'
$ws.Range('E44').Value = '  -5.32%  '
$ws.Range('D45').Value = '3.731'
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('D46').Value = '12.74'
$ws.Range('E46').Value = '  -6.35%  '
$ws.Range('D47').Value = '0.5479'
$ws.Range('E47').Value = '  -5.93%  '
$ws.Range('D48').Value = '116.23'
$ws.Range('E48').Value = '  -4.67%  '
$ws.Range('D49').Value = '1.847'
$ws.Range('E49').Value = '  -6.89%  '
$ws.Range('D50').Value = '1.122'
$ws.Range('E50').Value = '  -4.23%  '
$ws.Range('D51').Value = '0.06770'
$ws.Range('E51').Value = '  -3.39%  '

# Restore column D to the default "Normal" style so no stray number-format
# style is left behind on cells (matches original workbook formatting).
$dRange.Style = "Normal"
